$d = $word.ActiveDocument

# --- Change 1: "<caption>Bastard ball</caption>" -> "<caption><fr>Bastarde</fr> ball</caption>" ---
# This text lives across two separate runs within one paragraph:
#   run A (Courier New / blue, xml-tag style): "<caption>"
#   run B (plain style):                       "Bastard ball"
# We must edit each run's text individually (via a narrowed Range) so that
# Word does not merge the differently-formatted runs into a single run.

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*<caption>Bastard ball</caption>*") {
        $targetPara = $para
        break
    }
}

if ($targetPara -ne $null) {
    $pStart = $targetPara.Range.Start

    # run A: the literal "<caption>" text (9 characters) -> "<caption><fr>"
    $runA = $d.Range($pStart, $pStart + 9)
    $runA.Find.Execute("<caption>", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "<caption><fr>", 2)

    # run B: immediately follows run A and originally reads "Bastard ball" (12 chars)
    $runBStart = $pStart + 13
    $runB = $d.Range($runBStart, $runBStart + 12)
    $runB.Find.Execute("Bastard ball", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "Bastarde</fr> ball", 2)
}

# --- Change 2: "Moyen ball" -> "Average ball" ---
# This text is fully contained within a single run and is unique in the
# document, so a normal document-wide Find/Replace is safe here.
$d.Content.Find.Execute("Moyen ball", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Average ball", 2)
